$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing row 2 (FAPs -> Ccl28 -> Ackr2 -> FAPs) down to row 3,
# since that data is preserved unchanged in the new layout.
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(3, $c).Value = $ws.Cells.Item(2, $c).Value()
}

# Now overwrite row 2 with the new "ECs -> Ccl28 -> Ackr2 -> FAPs" edge values.
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl28"
$ws.Range("C2").Value = "Ackr2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.1014825
$ws.Range("H2").Value = 0.202965
$ws.Range("I2").Value = 0.5602924329074438
$ws.Range("J2").Value = 0.4593111407816425
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2786473333333334
$ws.Range("N2").Value = 0.8359420000000001
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.028277828005
$ws.Range("R2").Value = 0.16966696803
$ws.Range("S2").Value = 0.5602924329074438
$ws.Range("T2").Value = 0.4593111407816425

Write-Output "done"
